$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 27 (row 39): MI vs CSK results - set the raw score inputs
$ws.Range("E39").Value = 100
$ws.Range("H39").Value = 0
$ws.Range("K39").Value = 80
$ws.Range("N39").Value = 50
$ws.Range("Q39").Value = 40
$ws.Range("T39").Value = 60
$ws.Range("W39").Value = 30
$ws.Range("Z39").Value = 70
$ws.Range("AC39").Value = 20

$wb.Save()
